$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1/J1 - copy the style from the neighboring header cell (H1)
# so the new cells reuse the existing bold/border/centered style, then set
# their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-6 for columns I (I0) and J (IF)
$ws.Range("I2").Value = 6
$ws.Range("J2").Value = 7

$ws.Range("I3").Value = 11
$ws.Range("J3").Value = 11

$ws.Range("I4").Value = 9
$ws.Range("J4").Value = 9

$ws.Range("I5").Value = 9
$ws.Range("J5").Value = 9

$ws.Range("I6").Value = 8
$ws.Range("J6").Value = 8
